$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: E2 formula should subtract D2 (rent) as well as C2 (taxes)
$ws.Range("E2").Formula = "=-2062.95-C2-D2"

# Update the sheet view selection (no longer scrolled to G1, selection on D10)
$ws.Range("D10").Select()
